$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Insert a new row at row 58, shifting existing rows 58:127 down to 59:128
$ws.Rows.Item(58).Insert()

# Populate the new row 58 with the new record
$ws.Cells.Item(58, 1).Value = 9
$ws.Cells.Item(58, 2).Value = "Vega Central Mapocho de Santiago"
$ws.Cells.Item(58, 3).Value = "Metropolitana"
$ws.Cells.Item(58, 4).Value = 44966
$ws.Cells.Item(58, 5).Value = 13
$ws.Cells.Item(58, 6).Value = "Fruta"
$ws.Cells.Item(58, 7).Value = 100101
$ws.Cells.Item(58, 8).Value = "Berries"
$ws.Cells.Item(58, 9).Value = 100101004
$ws.Cells.Item(58, 10).Value = "Frambuesa"
$ws.Cells.Item(58, 11).Value = "Sin especificar"
$ws.Cells.Item(58, 12).Value = "Primera"
$ws.Cells.Item(58, 13).Value = 350
$ws.Cells.Item(58, 14).Value = 7000
$ws.Cells.Item(58, 15).Value = 7000
$ws.Cells.Item(58, 16).Value = 7000
$ws.Cells.Item(58, 17).Value = "`$/bandeja 2 kilos"
$ws.Cells.Item(58, 18).Value = "Provincia de Curicó"
$ws.Cells.Item(58, 19).Value = 3500
$ws.Cells.Item(58, 20).Value = 2

# Apply the same number format (date style) as the D column uses elsewhere
$ws.Cells.Item(58, 4).NumberFormat = $ws.Cells.Item(59, 4).NumberFormat
